# Beach tournament workbook update:
#  - Team on Anmeldung row 17 (seed 16) becomes a "bye" entry (both players,
#    both team slots) instead of "Linus Dürst" / "Marwin Dürst".
#  - Match #1 (Resultate row 2) result entered: 2 : 0.
#  - Match #9 (Resultate row 10) result entered: 0 : 2.
#  - Match sheet court for match #4 (row 5) set to court 1.
#  - Selections / active sheet updated to reflect where the user was working.

$wb = $excel.ActiveWorkbook

$anmeldung = $wb.Worksheets.Item("Anmeldung")
$resultate = $wb.Worksheets.Item("Resultate")
$match     = $wb.Worksheets.Item("Match")

# --- Anmeldung: seed 16 (row 17) becomes a bye ---
$anmeldung.Range("B17").Value = "bye"
$anmeldung.Range("C17").Value = "bye"
$anmeldung.Range("E17").Value = "bye"
$anmeldung.Range("F17").Value = "bye"

# --- Resultate: enter match results ---
$resultate.Range("H2").Value = 2
$resultate.Range("J2").Value = 0

$resultate.Range("H10").Value = 0
$resultate.Range("J10").Value = 2

# --- Match: set the court number for match #4 ---
$match.Range("C5").Value = 1

# --- Restore selections / active sheet to match where the user left off ---
$resultate.Activate()
$resultate.Range("H10").Select()

$match.Activate()
$match.Range("C6").Select()

$anmeldung.Activate()
$anmeldung.Range("G2").Select()
